$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header P1 from "implemented_stay_at_home" to "imposed_intervention"
$ws.Range("P1").Value = "imposed_intervention"

# 1b. Q2 used to describe Q1 (high_chance_of_death) as "tbc (boolean)" - it is now
#     resolved/defined as a plain boolean like the other y/n columns.
$ws.Range("Q2").Value = "y (boolean)"

# 2. Remove the yellow highlight fill from Q1 (keep bold font), so it matches the
#    plain bold style used by the other header cells (N1, O1, P1, L1, M1).
$ws.Range("N1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Add a new data row (row 4) with "y" values across most columns
$ws.Range("A4:M4").Value = "y"
$ws.Range("P4:Q4").Value = "y"

# 4. Autofit the columns so their widths reflect the new content (bestFit)
$ws.Range("A:P").EntireColumn.AutoFit()

# 5. Move the active selection to I10
$ws.Range("I10").Select()
